$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the "Ciudad" labels for rows 56 and 57 (Lanzarote <-> La Palma)
# and their corresponding "Muertes" (E column) counts, so that the
# row that used to read "Lanzarote"/3 now reads "La Palma"/4, and the
# row that used to read "La Palma"/4 now reads "Lanzarote"/3.
$ws.Cells.Item(56, 1).Value = "La Palma"
$ws.Cells.Item(57, 1).Value = "Lanzarote"

$ws.Cells.Item(56, 5).Value = 4
$ws.Cells.Item(57, 5).Value = 3

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 5 de Abril de 2020 a las 19:52"
